# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Halicarnassus_Profits sheets (columns H-N: price & profit figures).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H28").Value = 13855.889
$ws.Range("I28").Value = 1602.5
$ws.Range("K28").Value = 1602.5
$ws.Range("M28").Value = -1117.5
$ws.Range("H29").Value = 997.5
$ws.Range("J29").Value = 997.5
$ws.Range("L29").Value = 2992.5
$ws.Range("N29").Value = -3554.5
$ws.Range("H33").Value = 127.55556
$ws.Range("I33").Value = 127.55556
$ws.Range("K33").Value = 127.55556
$ws.Range("M33").Value = 101.44444
$ws.Range("H40").Value = 5189.722
$ws.Range("I40").Value = 3876.1667
$ws.Range("K40").Value = 3876.1667
$ws.Range("M40").Value = -3701.1667
$ws.Range("H138").Value = 2057.7307
$ws.Range("I138").Value = 987.75
$ws.Range("J138").Value = 3769.7
$ws.Range("K138").Value = 2963.25
$ws.Range("L138").Value = 11309.1
$ws.Range("M138").Value = 2176.75
$ws.Range("N138").Value = -21589.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4681.3335
$ws.Range("I2").Value = 2799.111
$ws.Range("K2").Value = 2799.111
$ws.Range("M2").Value = -2686.111
$ws.Range("H32").Value = 2506381.2
$ws.Range("I32").Value = 3215.0645
$ws.Range("K32").Value = 3215.0645
$ws.Range("M32").Value = -2928.0645
$ws.Range("H43").Value = 10033333
$ws.Range("I43").Value = 15000000
$ws.Range("K43").Value = 15000000
$ws.Range("M43").Value = -14999687
$ws.Range("H45").Value = 3567.2727
$ws.Range("I45").Value = 2498.5715
$ws.Range("K45").Value = 2498.5715
$ws.Range("M45").Value = -2121.5715
$ws.Range("H61").Value = 6750
$ws.Range("I61").Value = 6750
$ws.Range("K61").Value = 6750
$ws.Range("M61").Value = -6538
$ws.Range("H116").Value = 4681.3335
$ws.Range("I116").Value = 2799.111
$ws.Range("K116").Value = 2799.111
$ws.Range("M116").Value = -505.1109999999999
$ws.Range("H122").Value = 3967.5715
$ws.Range("I122").Value = 3630.1667
$ws.Range("K122").Value = 10890.5001
$ws.Range("M122").Value = -8440.500100000001
$ws.Range("H136").Value = 6750
$ws.Range("I136").Value = 6750
$ws.Range("K136").Value = 20250
$ws.Range("M136").Value = -17700
$ws.Range("H139").Value = 50714.5
$ws.Range("J139").Value = 50714.5
$ws.Range("L139").Value = 50714.5
$ws.Range("N139").Value = -60994.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4681.3335
$ws.Range("I3").Value = 2799.111
$ws.Range("K3").Value = 2799.111
$ws.Range("M3").Value = -2685.111
$ws.Range("H94").Value = 900
$ws.Range("I94").Value = 825
$ws.Range("K94").Value = 825
$ws.Range("M94").Value = -374
$ws.Range("H95").Value = 9812
$ws.Range("J95").Value = 9812
$ws.Range("L95").Value = 9812
$ws.Range("N95").Value = -15304
$ws.Range("H105").Value = 2769.7778
$ws.Range("I105").Value = 1864.625
$ws.Range("K105").Value = 1864.625
$ws.Range("M105").Value = -117.625
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -10524
$ws.Range("H134").Value = 2692.8635
$ws.Range("I134").Value = 2381.2632
$ws.Range("K134").Value = 7143.7896
$ws.Range("M134").Value = -4608.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1174.75
$ws.Range("I16").Value = 1166.3334
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1166.3334
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -879.3334
$ws.Range("N16").Value = -1774
$ws.Range("H22").Value = 1386.6923
$ws.Range("J22").Value = 1554.7142
$ws.Range("L22").Value = 1554.7142
$ws.Range("N22").Value = -2254.7142
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384
$ws.Range("H94").Value = 3004.0557
$ws.Range("I94").Value = 2404.3333
$ws.Range("K94").Value = 2404.3333
$ws.Range("M94").Value = -1953.3333
$ws.Range("H113").Value = 1174.75
$ws.Range("I113").Value = 1166.3334
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1166.3334
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1003.6666
$ws.Range("N113").Value = -5540
$ws.Range("H125").Value = 67925
$ws.Range("J125").Value = 67925
$ws.Range("L125").Value = 67925
$ws.Range("N125").Value = -72845
$ws.Range("H132").Value = 4184.923
$ws.Range("I132").Value = 3526.087
$ws.Range("J132").Value = 5132
$ws.Range("K132").Value = 10578.261
$ws.Range("L132").Value = 15396
$ws.Range("M132").Value = -8048.261
$ws.Range("N132").Value = -20456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 278.8
$ws.Range("I8").Value = 278.8
$ws.Range("K8").Value = 836.4000000000001
$ws.Range("M8").Value = -697.4000000000001
$ws.Range("H50").Value = 298.75
$ws.Range("I50").Value = 47.5
$ws.Range("K50").Value = 142.5
$ws.Range("M50").Value = 338.5
$ws.Range("H53").Value = 298.75
$ws.Range("I53").Value = 47.5
$ws.Range("K53").Value = 142.5
$ws.Range("M53").Value = 338.5
$ws.Range("H132").Value = 2132.3333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2132.3333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19190.9997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -24250.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H80").Value = 2900
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2900
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1902
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2900
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 14500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -9508
$ws.Range("N83").ClearContents()
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492
$ws.Range("H113").Value = 8000.143
$ws.Range("I113").Value = 5334.1665
$ws.Range("K113").Value = 5334.1665
$ws.Range("M113").Value = -3164.1665
$ws.Range("H122").Value = 386260.53
$ws.Range("I122").Value = 418033.1
$ws.Range("J122").Value = 4990
$ws.Range("K122").Value = 1254099.3
$ws.Range("L122").Value = 14970
$ws.Range("M122").Value = -1251649.3
$ws.Range("N122").Value = -19870
$ws.Range("H132").Value = 36394.03
$ws.Range("I132").Value = 44660.64
$ws.Range("J132").Value = 6870.4287
$ws.Range("K132").Value = 133981.92
$ws.Range("L132").Value = 20611.2861
$ws.Range("M132").Value = -131451.92
$ws.Range("N132").Value = -25671.2861
$ws.Range("H136").Value = 10000
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 350
$ws.Range("I2").Value = 350
$ws.Range("K2").Value = 350
$ws.Range("M2").Value = -238
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H82").Value = 4690.909
$ws.Range("I82").Value = 2924.5
$ws.Range("J82").Value = 5700.2856
$ws.Range("K82").Value = 2924.5
$ws.Range("L82").Value = 5700.2856
$ws.Range("M82").Value = -2563.5
$ws.Range("N82").Value = -6422.2856
$ws.Range("H85").Value = 4690.909
$ws.Range("I85").Value = 2924.5
$ws.Range("J85").Value = 5700.2856
$ws.Range("K85").Value = 2924.5
$ws.Range("L85").Value = 5700.2856
$ws.Range("M85").Value = -1676.5
$ws.Range("N85").Value = -8196.285599999999
$ws.Range("H122").Value = 5003.3
$ws.Range("I122").Value = 5010.5
$ws.Range("J122").Value = 4974.5
$ws.Range("K122").Value = 15031.5
$ws.Range("L122").Value = 14923.5
$ws.Range("M122").Value = -12581.5
$ws.Range("N122").Value = -19823.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9752
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9142
$ws.Range("H81").Value = 712.5
$ws.Range("I81").Value = 712.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1425
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -364
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 712.5
$ws.Range("I84").Value = 712.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7125
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1821
$ws.Range("N84").ClearContents()
